# Generate Report for Handoff
#
# Updates the "b.md" row across the Overview / zh-cn / de-de sheets to
# reflect a brand new handoff (new xlf file, new status, new timestamps).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" file.
# Status columns (zh-cn, de-de) move from "Handed back: in sync with en-US"
# to "Ready for handoff", and the Latest Handoff Date is refreshed.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-24 02:37:22"

# ---------------------------------------------------------------------
# zh-cn detail sheet: row 3 is the "b.md" file.
# Status -> "Ready for handoff"
# Latest Handoff File -> new xlf file name (hyperlink text also updated)
# Latest Handoff Datetime -> new timestamp
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-24 02:37:18"

foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de detail sheet: row 3 is the "b.md" file.
# Status -> "Ready for handoff"
# Latest Handoff File -> new xlf file name (hyperlink text also updated)
# Latest Handoff Datetime -> new timestamp
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-24 02:37:22"

foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
